# "Added the Widgets Page And TestCases"
#
# Replace the old auto/color/multi-auto test data (rows 1-12, cols A-I)
# with the new, smaller "countries / countryNames" test data (rows 1-4,
# cols A-B), rename the sheet/project, and keep the autofilter + the
# hidden _FilterDatabase name in sync with the new, smaller range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data values --------------------------------------------------
# Write order controls shared-string allocation order, so write the
# first occurrence of each distinct string in the same order the
# authored workbook uses: countries, "United States,...", countryNames,
# "New Zealand,...".
$ws.Range("A1").Value = "countries"
$ws.Range("B2").Value = "United States,Panama,Zimbabwe"
$ws.Range("B1").Value = "countryNames"
$ws.Range("B3").Value = "New Zealand,Argentina,Puerto Rico"
$ws.Range("A2").Value = "countries"
$ws.Range("A3").Value = "countries"
$ws.Range("A4").Value = "countries"
$ws.Range("B4").Value = "United States,Panama,Zimbabwe"

# --- drop the old rows/columns that are no longer part of the sheet ---
$ws.Rows("5:12").Delete()
$ws.Columns("C:I").Delete()

# --- column width for column B (column A keeps its existing width) -----
$ws.Columns.Item(2).ColumnWidth = 36.86

# --- rename sheet/tab to match the new project name --------------------
$ws.Name = "in.automationtest"

# --- re-apply the autofilter over the smaller A1:B1 range --------------
$ws.AutoFilterMode = $false
$ws.Range("A1:B1").AutoFilter()

# keep the hidden _xlnm._FilterDatabase name in sync with the new range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=in.automationtest!`$A`$1:`$B`$1"
    }
}

# --- selection lands on B3, matching the authored sheetView ------------
$ws.Range("B3").Select()
